$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsSchedule.Range("B2").Value = 46043.16666666666
$wsSchedule.Range("C2").Value = 4
$wsSchedule.Range("D2").Value = 15.12
$wsSchedule.Range("E2").Value = 506.74358475
$wsSchedule.Range("F2").Value = 33.51478735119048
$wsSchedule.Range("A4").Value = 46043.875
$wsSchedule.Range("B4").Value = 46044.08333333334
$wsSchedule.Range("E4").Value = 667.8974639999999
$wsSchedule.Range("F4").Value = 35.33849015873015
$wsSchedule.Range("A5").Value = 46044.25
$wsSchedule.Range("C5").Value = 10
$wsSchedule.Range("D5").Value = 37.8
$wsSchedule.Range("E5").Value = 54.28589400000001
$wsSchedule.Range("F5").Value = 1.436134761904762
$wsDetailed.Range("E10").Value = "OFF"
$wsDetailed.Range("B37").Value = 48.61802
$wsDetailed.Range("B38").Value = 55.37363
$wsDetailed.Range("B39").Value = 61.05003
$wsDetailed.Range("C39").Value = "historical"
$wsDetailed.Range("B40").Value = 79.95
$wsDetailed.Range("C40").Value = "historical"
$wsDetailed.Range("B41").Value = 79.95
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 73.37
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("B43").Value = 65
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("B44").Value = 57.31
$wsDetailed.Range("C44").Value = "historical"
$wsDetailed.Range("E44").Value = "ON"
$wsDetailed.Range("B45").Value = 59.66383
$wsDetailed.Range("C45").Value = "historical"
$wsDetailed.Range("B46").Value = 57.09
$wsDetailed.Range("C46").Value = "historical"
$wsDetailed.Range("B47").Value = 90.16674999999999
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("B48").Value = 76.08317
$wsDetailed.Range("C48").Value = "historical"
$wsDetailed.Range("B49").Value = 73.2
$wsDetailed.Range("B50").Value = 65.84798000000001
$wsDetailed.Range("B51").Value = 73.2
$wsDetailed.Range("B52").Value = 66.16679999999999
$wsDetailed.Range("B53").Value = 66.29451
$wsDetailed.Range("B54").Value = 65.19963
$wsDetailed.Range("E54").Value = "OFF"
$wsDetailed.Range("B55").Value = 66.09792
$wsDetailed.Range("B56").Value = 67.38155
$wsDetailed.Range("B57").Value = 73.2
$wsDetailed.Range("B58").Value = 73.2
$wsDetailed.Range("B59").Value = 73.2
$wsDetailed.Range("B60").Value = 66.36660999999999
$wsDetailed.Range("B61").Value = 78
$wsDetailed.Range("E62").Value = "ON"
$wsDetailed.Range("B63").Value = 57.31
$wsDetailed.Range("B64").Value = 35.88
$wsDetailed.Range("B66").Value = -5.74313
$wsDetailed.Range("B67").Value = -6.44675
$wsDetailed.Range("B68").Value = -8.05181
$wsDetailed.Range("B69").Value = -7.78632
$wsDetailed.Range("B70").Value = -7.66245
$wsDetailed.Range("B71").Value = -7.9504
$wsDetailed.Range("B72").Value = -8.920400000000001
$wsDetailed.Range("B73").Value = -5.1817
$wsDetailed.Range("B74").Value = -7.39026
$wsDetailed.Range("B75").Value = -6.31903
$wsDetailed.Range("B76").Value = -7.86005
$wsDetailed.Range("B77").Value = -5.91747
$wsDetailed.Range("B78").Value = -5.3165
$wsDetailed.Range("B79").Value = -0.86589
$wsDetailed.Range("B81").Value = -12.01
$wsDetailed.Range("B82").Value = -9.71002
$wsDetailed.Range("B83").Value = -10
$wsDetailed.Range("B84").Value = -11.25715
$wsDetailed.Range("B85").Value = -8.19045
$wsDetailed.Range("B86").Value = -1.63851
$wsDetailed.Range("B87").Value = 7.77782
$wsDetailed.Range("B88").Value = 18.58582
$wsDetailed.Range("B90").Value = 55.33037
$wsDetailed.Range("B91").Value = 52.96163
$wsDetailed.Range("B93").Value = 57.03541
$wsDetailed.Range("B94").Value = 42.9873
$wsDetailed.Range("B96").Value = 56.98
$wsDetailed.Range("B97").Value = 48.64303
